# Graham Potter sheet: add a new "Team" column (AF) indicating the club
# he was managing for each match row. For this workbook every row involves
# Swansea (his club during this spell), so the new column is "Swansea"
# for every data row, with header "Team" in AF1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell AF1 = "Team", styled like the other header cells (bold,
# centered, bordered) by copying the format from the preceding header
# cell AE1.
$ws.Range("AF1").Value = "Team"
$ws.Range("AE1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)

# Data rows 2-139: AF = "Swansea"
$ws.Range("AF2:AF139").Value = "Swansea"
